$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows whose Target cluster (column D) is "ECs"
# Delete from the bottom up so row indices of earlier rows are unaffected
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(2).Delete()

# Refresh remaining 6 data rows (now rows 2-7) with the updated TPM-derived values

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lrpap1"
$ws.Range("C2").Value = "Lrp8"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.573811
$ws.Range("H2").Value = 13.721433
$ws.Range("I2").Value = 0.1659009079913533
$ws.Range("J2").Value = 0.1659009079913533
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1124773333333333
$ws.Range("N2").Value = 0.337432
$ws.Range("O2").Value = 0.7871437602495106
$ws.Range("P2").Value = 0.7871437602495107
$ws.Range("Q2").Value = 0.5144500644506667
$ws.Range("R2").Value = 4.630050580055999
$ws.Range("S2").Value = 0.1305878645451219
$ws.Range("T2").Value = 0.1305878645451219

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lrpap1"
$ws.Range("C3").Value = "Lrp8"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.573811
$ws.Range("H3").Value = 13.721433
$ws.Range("I3").Value = 0.1659009079913533
$ws.Range("J3").Value = 0.1659009079913533
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03041566666666666
$ws.Range("N3").Value = 0.09124699999999999
$ws.Range("O3").Value = 0.2128562397504893
$ws.Range("P3").Value = 0.2128562397504893
$ws.Range("Q3").Value = 0.1391155107723333
$ws.Range("R3").Value = 1.252039596951
$ws.Range("S3").Value = 0.03531304344623137
$ws.Range("T3").Value = 0.03531304344623137

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lrpap1"
$ws.Range("C4").Value = "Lrp8"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 14.67485766666667
$ws.Range("H4").Value = 44.024573
$ws.Range("I4").Value = 0.5322852674812913
$ws.Range("J4").Value = 0.5322852674812913
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1124773333333333
$ws.Range("N4").Value = 0.337432
$ws.Range("O4").Value = 0.7871437602495106
$ws.Range("P4").Value = 0.7871437602495107
$ws.Range("Q4").Value = 1.650588857392889
$ws.Range("R4").Value = 14.855299716536
$ws.Range("S4").Value = 0.4189850269706402
$ws.Range("T4").Value = 0.4189850269706403

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lrpap1"
$ws.Range("C5").Value = "Lrp8"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 14.67485766666667
$ws.Range("H5").Value = 44.024573
$ws.Range("I5").Value = 0.5322852674812913
$ws.Range("J5").Value = 0.5322852674812913
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.03041566666666666
$ws.Range("N5").Value = 0.09124699999999999
$ws.Range("O5").Value = 0.2128562397504893
$ws.Range("P5").Value = 0.2128562397504893
$ws.Range("Q5").Value = 0.4463455791701111
$ws.Range("R5").Value = 4.017110212531001
$ws.Range("S5").Value = 0.1133002405106511
$ws.Range("T5").Value = 0.1133002405106511

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Lrpap1"
$ws.Range("C6").Value = "Lrp8"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 8.320867
$ws.Range("H6").Value = 24.962601
$ws.Range("I6").Value = 0.3018138245273554
$ws.Range("J6").Value = 0.3018138245273554
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1124773333333333
$ws.Range("N6").Value = 0.337432
$ws.Range("O6").Value = 0.7871437602495106
$ws.Range("P6").Value = 0.7871437602495107
$ws.Range("Q6").Value = 0.9359089311813332
$ws.Range("R6").Value = 8.423180380631999
$ws.Range("S6").Value = 0.2375708687337485
$ws.Range("T6").Value = 0.2375708687337485

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Lrpap1"
$ws.Range("C7").Value = "Lrp8"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 8.320867
$ws.Range("H7").Value = 24.962601
$ws.Range("I7").Value = 0.3018138245273554
$ws.Range("J7").Value = 0.3018138245273554
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.03041566666666666
$ws.Range("N7").Value = 0.09124699999999999
$ws.Range("O7").Value = 0.2128562397504893
$ws.Range("P7").Value = 0.2128562397504893
$ws.Range("Q7").Value = 0.2530847170496666
$ws.Range("R7").Value = 2.277762453447
$ws.Range("S7").Value = 0.06424295579360687
$ws.Range("T7").Value = 0.06424295579360687
